$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.118.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.787.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '344.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '115.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0854'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.132'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.227.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.811.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.886'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.037.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.83%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.23%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("E25").Value = '  +6.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("E28").Value = '  -0.79%  '

$ws.Range("E29").Value = '  +0.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.141'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.35%  '

$ws.Range("E33").Value = '  +1.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0824'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0413'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +18.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.37%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +22.34%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.73%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.116'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '127.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.02%  '

$ws.Range("E45").Value = '  +0.27%  '

$ws.Range("E46").Value = '  -2.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.070.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.87%  '

$ws.Range("E48").Value = '  +1.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.901'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +13.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.86%  '
